# Adds new rows of synthesis-center/network data to the data map worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New records to append starting at row 20.
# Columns: A=Name, B=Abbreviation, C=City, D=State, E=Country, F=Abb_Country, G=Lat, H=Long, I=Active
$rows = @(
    @{ A = "Intergovernmental Science-Policy Platform on Biodiversity and Ecosystem Services"; B = "IPBES";   C = "Bonn";       D = "NA";                    E = "Germany";                   F = "GER"; G = 50.71848;           H = 7.1254629999999999;   I = "Yes" },
    @{ A = "Science for Nature and People Partnership";                                         B = "SNAPP";   C = "Arlington";  D = "Virginia";              E = "United States of America";  F = "USA"; G = 38.877848;          H = -77.089731999999998;  I = "Yes" },
    @{ A = "EU Knowledge and Learning Mechanism on Biodiversity and Ecosystem Services";         B = "EKLIPSE"; C = "Leipzig";    D = "NA";                    E = "Germany";                   F = "GER"; G = 51.351638000000001; H = 12.430899999999999;   I = "Yes" },
    @{ A = "Lenfest Ocean Programme at The Pew Charitable Trusts";                               B = "Lenfest"; C = "Washington"; D = "District of Columbia"; E = "United States of America";  F = "USA"; G = 38.897542999999999; H = -77.026568999999995;  I = "Yes" }
)

# Insert the new rows right after the existing data (row 19) so that each new
# row inherits the same cell formatting/style pattern the author used for the
# previous rows (columns C, D, E, F, I styled; A, B, G, H left on default style).
$startRow = 20
for ($i = 0; $i -lt $rows.Count; $i++) {
    $ws.Rows.Item($startRow).Insert()
}

# Row 20 was typed Abbreviation/City first, then Name (matches the shared-string
# insertion order in the source file: IPBES, Bonn, then the long Name string).
# Rows 21-23 were typed left-to-right starting with Name.
$rec = $rows[0]
$r = $startRow
$ws.Cells.Item($r, 2).Value = $rec.B
$ws.Cells.Item($r, 3).Value = $rec.C
$ws.Cells.Item($r, 1).Value = $rec.A
$ws.Cells.Item($r, 4).Value = $rec.D
$ws.Cells.Item($r, 5).Value = $rec.E
$ws.Cells.Item($r, 6).Value = $rec.F
$ws.Cells.Item($r, 7).Value = $rec.G
$ws.Cells.Item($r, 8).Value = $rec.H
$ws.Cells.Item($r, 9).Value = $rec.I

for ($i = 1; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rec = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $rec.A
    $ws.Cells.Item($r, 2).Value = $rec.B
    $ws.Cells.Item($r, 3).Value = $rec.C
    $ws.Cells.Item($r, 4).Value = $rec.D
    $ws.Cells.Item($r, 5).Value = $rec.E
    $ws.Cells.Item($r, 6).Value = $rec.F
    $ws.Cells.Item($r, 7).Value = $rec.G
    $ws.Cells.Item($r, 8).Value = $rec.H
    $ws.Cells.Item($r, 9).Value = $rec.I
}

# D23 was left unstyled by the author (unlike D20:D22), matching the data entry
# inconsistency visible in the source workbook.
$ws.Cells.Item(23, 4).ClearFormats()

# Move the active selection to A24, as left by the author after data entry.
$ws.Range("A24").Select()

$wb.Save()
